$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data rows (2-13) have had their
# Fecha(D) / Volumen(M) / Precio minimo(N) / Precio maximo(O) /
# Precio promedio ponderado(P) / Precio $/Kg(S) values re-shuffled
# between rows. Row 7 stays the same.
#
# Mapping: new row <- source (old) row
#   2  <- 6
#   3  <- 11
#   4  <- 12
#   5  <- 3
#   6  <- 8
#   7  <- 7   (unchanged)
#   8  <- 5
#   9  <- 13
#   10 <- 9
#   11 <- 10
#   12 <- 4
#   13 <- 2

# Capture original values (D, M, N, O, P, S) for rows 2..13 before overwriting.
$orig = @{}
for ($r = 2; $r -le 13; $r++) {
    $orig[$r] = @(
        $ws.Range("D$r").Value2,
        $ws.Range("M$r").Value2,
        $ws.Range("N$r").Value2,
        $ws.Range("O$r").Value2,
        $ws.Range("P$r").Value2,
        $ws.Range("S$r").Value2
    )
}

$mapping = @{
    2  = 6
    3  = 11
    4  = 12
    5  = 3
    6  = 8
    7  = 7
    8  = 5
    9  = 13
    10 = 9
    11 = 10
    12 = 4
    13 = 2
}

foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    $vals = $orig[$src]
    $ws.Range("D$r").Value2 = $vals[0]
    $ws.Range("M$r").Value2 = $vals[1]
    $ws.Range("N$r").Value2 = $vals[2]
    $ws.Range("O$r").Value2 = $vals[3]
    $ws.Range("P$r").Value2 = $vals[4]
    $ws.Range("S$r").Value2 = $vals[5]
}
